$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3834.5
$ws.Range("I76").Value = 3000.75
$ws.Range("J76").Value = 5502
$ws.Range("K76").Value = 3000.75
$ws.Range("L76").Value = 5502
$ws.Range("M76").Value = -2685.75
$ws.Range("N76").Value = -6132
$ws.Range("H79").Value = 3834.5
$ws.Range("I79").Value = 3000.75
$ws.Range("J79").Value = 5502
$ws.Range("K79").Value = 3000.75
$ws.Range("L79").Value = 5502
$ws.Range("M79").Value = -1908.75
$ws.Range("N79").Value = -7686
$ws.Range("H87").Value = 20328
$ws.Range("J87").Value = 20328
$ws.Range("L87").Value = 20328
$ws.Range("N87").Value = -22824
$ws.Range("H90").Value = 20328
$ws.Range("J90").Value = 20328
$ws.Range("L90").Value = 60984
$ws.Range("N90").Value = -73464
$ws.Range("H132").Value = 5420.9033
$ws.Range("I132").Value = 5297.8335
$ws.Range("J132").Value = 5842.857
$ws.Range("K132").Value = 15893.5005
$ws.Range("L132").Value = 17528.571
$ws.Range("M132").Value = -13363.5005
$ws.Range("N132").Value = -22588.571
$ws.Range("H138").Value = 2386.5356
$ws.Range("I138").Value = 1199.75
$ws.Range("J138").Value = 3573.3215
$ws.Range("K138").Value = 3599.25
$ws.Range("L138").Value = 10719.9645
$ws.Range("M138").Value = 1540.75
$ws.Range("N138").Value = -20999.9645

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6107.3857
$ws.Range("I32").Value = 5285.302
$ws.Range("K32").Value = 5285.302
$ws.Range("M32").Value = -4998.302
$ws.Range("H45").Value = 1264.4
$ws.Range("I45").Value = 1264.4
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1264.4
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -887.4000000000001
$ws.Range("H88").Value = 2475.7334
$ws.Range("I88").Value = 1472.6666
$ws.Range("J88").Value = 3144.4443
$ws.Range("K88").Value = 1472.6666
$ws.Range("L88").Value = 3144.4443
$ws.Range("M88").Value = -1066.6666
$ws.Range("N88").Value = -3956.4443
$ws.Range("H91").Value = 2475.7334
$ws.Range("I91").Value = 1472.6666
$ws.Range("J91").Value = 3144.4443
$ws.Range("K91").Value = 1472.6666
$ws.Range("L91").Value = 3144.4443
$ws.Range("M91").Value = -68.66660000000002
$ws.Range("N91").Value = -5952.4443
$ws.Range("H139").Value = 65539.25
$ws.Range("J139").Value = 65539.25
$ws.Range("L139").Value = 65539.25
$ws.Range("N139").Value = -75819.25
$ws.Range("N45").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1207.2174
$ws.Range("I20").Value = 1286.5405
$ws.Range("J20").Value = 881.1111
$ws.Range("K20").Value = 1286.5405
$ws.Range("L20").Value = 881.1111
$ws.Range("M20").Value = -1039.5405
$ws.Range("N20").Value = -1375.1111
$ws.Range("H86").Value = 2698.7144
$ws.Range("I86").Value = 2752.4614
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 2752.4614
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -1629.4614
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 2698.7144
$ws.Range("I89").Value = 2752.4614
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 13762.307
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -8146.307000000001
$ws.Range("N89").Value = -21232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4563.1943
$ws.Range("I31").Value = 4875
$ws.Range("J31").Value = 3627.7778
$ws.Range("K31").Value = 4875
$ws.Range("L31").Value = 3627.7778
$ws.Range("M31").Value = -4580
$ws.Range("N31").Value = -4217.7778
$ws.Range("H34").Value = 4563.1943
$ws.Range("I34").Value = 4875
$ws.Range("J34").Value = 3627.7778
$ws.Range("K34").Value = 4875
$ws.Range("L34").Value = 3627.7778
$ws.Range("M34").Value = -4673
$ws.Range("N34").Value = -4031.7778
$ws.Range("H62").Value = 2415.1428
$ws.Range("I62").Value = 2400
$ws.Range("K62").Value = 2400
$ws.Range("M62").Value = -1776
$ws.Range("H65").Value = 2415.1428
$ws.Range("I65").Value = 2400
$ws.Range("K65").Value = 12000
$ws.Range("M65").Value = -8880
$ws.Range("H107").Value = 796.1739
$ws.Range("I107").Value = 1349.2
$ws.Range("J107").Value = 370.76923
$ws.Range("K107").Value = 1349.2
$ws.Range("L107").Value = 370.76923
$ws.Range("M107").Value = 570.8
$ws.Range("N107").Value = -4210.76923
$ws.Range("H127").Value = 40637.375
$ws.Range("J127").Value = 40637.375
$ws.Range("L127").Value = 40637.375
$ws.Range("N127").Value = -50557.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4000
$ws.Range("I56").Value = 4000
$ws.Range("K56").Value = 4000
$ws.Range("M56").Value = -3470
$ws.Range("H96").Value = 4013.8333
$ws.Range("J96").Value = 4013.8333
$ws.Range("L96").Value = 12041.4999
$ws.Range("N96").Value = -16159.4999
$ws.Range("H113").Value = 1392.1765
$ws.Range("I113").Value = 2609.5715
$ws.Range("J113").Value = 540
$ws.Range("K113").Value = 7828.7145
$ws.Range("L113").Value = 1620
$ws.Range("M113").Value = -5658.7145
$ws.Range("N113").Value = -5960

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5181.5
$ws.Range("I70").Value = 5018
$ws.Range("J70").Value = 5279.6
$ws.Range("K70").Value = 5018
$ws.Range("L70").Value = 5279.6
$ws.Range("M70").Value = -4748
$ws.Range("N70").Value = -5819.6
$ws.Range("H73").Value = 5181.5
$ws.Range("I73").Value = 5018
$ws.Range("J73").Value = 5279.6
$ws.Range("K73").Value = 5018
$ws.Range("L73").Value = 5279.6
$ws.Range("M73").Value = -4082
$ws.Range("N73").Value = -7151.6
$ws.Range("H105").Value = 30000
$ws.Range("J105").Value = 30000
$ws.Range("L105").Value = 30000
$ws.Range("H126").Value = 92121.09
$ws.Range("I126").Value = 200938.4
$ws.Range("J126").Value = 1440
$ws.Range("K126").Value = 602815.2
$ws.Range("L126").Value = 4320
$ws.Range("M126").Value = -600345.2
$ws.Range("N126").Value = -9260
$ws.Range("N105").Value = -36988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("H136").Value = 9990.166999999999
$ws.Range("I136").Value = 12853.556
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 38560.66800000001
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -36010.66800000001
$ws.Range("N136").Value = -9300
$ws.Range("M107").ClearContents()
